$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at O, shifting existing O:V columns (and their
# formatting/widths) one place to the right, becoming P:W.
$null = $ws.Columns("O:O").Insert()

# Populate the new column's header with the new field name.
$ws.Cells.Item(1, 15).Value = "DOCREF4"

# Remove the now-stray last row (row 5), which only held a style-only cell.
$null = $ws.Rows("5:5").Delete()

# Update the active selection to match the target state.
$null = $ws.Range("I7").Select()
